$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.292742
$ws.Range("H2").Value = 0.8782260000000001
$ws.Range("I2").Value = 0.1005821958520865
$ws.Range("J2").Value = 0.1005821958520865
$ws.Range("M2").Value = 0.006825333333333333
$ws.Range("N2").Value = 0.020476
$ws.Range("O2").Value = 0.0506148883313352
$ws.Range("P2").Value = 0.0506148883313352
$ws.Range("Q2").Value = 0.001998061730666667
$ws.Range("R2").Value = 0.017982555576
$ws.Range("S2").Value = 0.005090956611173843
$ws.Range("T2").Value = 0.005090956611173845
# Row 3
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.292742
$ws.Range("H3").Value = 0.8782260000000001
$ws.Range("I3").Value = 0.1005821958520865
$ws.Range("J3").Value = 0.1005821958520865
$ws.Range("O3").Value = 0.3671161428271267
$ws.Range("P3").Value = 0.3671161428271267
$ws.Range("Q3").Value = 0.01449219271
$ws.Range("R3").Value = 0.13042973439
$ws.Range("S3").Value = 0.03692534777830061
$ws.Range("T3").Value = 0.03692534777830062
# Row 4
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.292742
$ws.Range("H4").Value = 0.8782260000000001
$ws.Range("I4").Value = 0.1005821958520865
$ws.Range("J4").Value = 0.1005821958520865
$ws.Range("M4").Value = 0.078518
$ws.Range("N4").Value = 0.235554
$ws.Range("O4").Value = 0.582268968841538
$ws.Range("P4").Value = 0.582268968841538
$ws.Range("Q4").Value = 0.022985516356
$ws.Range("R4").Value = 0.206869647204
$ws.Range("S4").Value = 0.05856589146261201
$ws.Range("T4").Value = 0.05856589146261203
# Row 5
$ws.Range("I5").Value = 0.8949043375045497
$ws.Range("J5").Value = 0.8949043375045498
$ws.Range("M5").Value = 0.006825333333333333
$ws.Range("N5").Value = 0.020476
$ws.Range("O5").Value = 0.0506148883313352
$ws.Range("P5").Value = 0.0506148883313352
$ws.Range("Q5").Value = 0.017777242724
$ws.Range("R5").Value = 0.159995184516
$ws.Range("S5").Value = 0.04529548311002029
$ws.Range("T5").Value = 0.0452954831100203
# Row 6
$ws.Range("I6").Value = 0.8949043375045497
$ws.Range("J6").Value = 0.8949043375045498
$ws.Range("O6").Value = 0.3671161428271267
$ws.Range("P6").Value = 0.3671161428271267
$ws.Range("S6").Value = 0.3285338285839355
$ws.Range("T6").Value = 0.3285338285839355
# Row 7
$ws.Range("I7").Value = 0.8949043375045497
$ws.Range("J7").Value = 0.8949043375045498
$ws.Range("M7").Value = 0.078518
$ws.Range("N7").Value = 0.235554
$ws.Range("O7").Value = 0.582268968841538
$ws.Range("P7").Value = 0.582268968841538
$ws.Range("Q7").Value = 0.204507747246
$ws.Range("R7").Value = 1.840569725214
$ws.Range("S7").Value = 0.5210750258105938
$ws.Range("T7").Value = 0.521075025810594
# Row 8
$ws.Range("E8").Value = 1.0
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.01313633333333333
$ws.Range("H8").Value = 0.039409
$ws.Range("I8").Value = 0.004513466643363867
$ws.Range("J8").Value = 0.004513466643363868
$ws.Range("M8").Value = 0.006825333333333333
$ws.Range("N8").Value = 0.020476
$ws.Range("O8").Value = 0.0506148883313352
$ws.Range("P8").Value = 0.0506148883313352
$ws.Range("Q8").Value = 0.00008965985377777778
$ws.Range("R8").Value = 0.0008069386840000001
$ws.Range("S8").Value = 0.0002284486101410684
$ws.Range("T8").Value = 0.0002284486101410685
# Row 9
$ws.Range("E9").Value = 1.0
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.01313633333333333
$ws.Range("H9").Value = 0.039409
$ws.Range("I9").Value = 0.004513466643363867
$ws.Range("J9").Value = 0.004513466643363868
$ws.Range("O9").Value = 0.3671161428271267
$ws.Range("P9").Value = 0.3671161428271267
$ws.Range("Q9").Value = 0.0006503141816666667
$ws.Range("R9").Value = 0.005852827635
$ws.Range("S9").Value = 0.001656966464890642
$ws.Range("T9").Value = 0.001656966464890642
# Row 10
$ws.Range("E10").Value = 1.0
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.01313633333333333
$ws.Range("H10").Value = 0.039409
$ws.Range("I10").Value = 0.004513466643363867
$ws.Range("J10").Value = 0.004513466643363868
$ws.Range("M10").Value = 0.078518
$ws.Range("N10").Value = 0.235554
$ws.Range("O10").Value = 0.582268968841538
$ws.Range("P10").Value = 0.582268968841538
$ws.Range("Q10").Value = 0.001031438620666667
$ws.Range("R10").Value = 0.009282947586
$ws.Range("S10").Value = 0.002628051568332157
$ws.Range("T10").Value = 0.002628051568332157

Write-Output "Applied all TPM updates"
